# Append the three new SQL-run log rows captured on Skylake-PC for the
# Latrunculi_ne test into the "test1" results sheet, then let the
# "compare1" sheet's lookup formulas pick the new rows up.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("test1")

$ws.Range("A6").Value = "2019-02-05 19:17:15"
$ws.Range("B6").Value = "Skylake-PC"
$ws.Range("C6").Value = "Latrunculi_ne"
$ws.Range("D6").Value = 3668
$ws.Range("E6").Value = 1.898438453674316

$ws.Range("A7").Value = "2019-02-05 19:20:23"
$ws.Range("B7").Value = "Skylake-PC"
$ws.Range("C7").Value = "Latrunculi_ne"
$ws.Range("D7").Value = 6735
$ws.Range("E7").Value = 2.942668676376343

$ws.Range("A8").Value = "2019-02-05 19:21:18"
$ws.Range("B8").Value = "Skylake-PC"
$ws.Range("C8").Value = "Latrunculi_ne"
$ws.Range("D8").Value = 4389
$ws.Range("E8").Value = 3.050406694412231

# Nudge the "compare1" lookup formulas so they pick up the new test1 rows
# (re-assigning a formula to itself forces recalculation of that range).
$ws2 = $wb.Worksheets.Item("compare1")
$rng = $ws2.Range("A6:E24")
$rng.Formula = $rng.Formula
